# Add files via upload
# Replace the old test whitelist emails with the real whitelist emails,
# append two more whitelisted senders, and strip the mailto: hyperlinks
# that Excel had auto-created for the email addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the auto-generated mailto hyperlinks (and their blue/underline
#     "Hyperlink" look) from the three original rows -------------------
$ws.Hyperlinks.Delete()
$ws.Range("A1").Style = "Normal"
$ws.Range("A3").Style = "Normal"

# --- Replace the placeholder addresses with the real whitelisted senders -
$ws.Range("A1").Value = "googleplay-noreply@google.com"
$ws.Range("A2").Value = "googledrive-noreply@google.com"
$ws.Range("A3").Value = "sc.singapore@sc.com"

# --- Append two more whitelisted senders in rows 4 and 5 -----------------
$ws.Range("A4").Value = "ibanking.alert@dbs.com"
$ws.Range("B4").Value = "whitelisted"

$ws.Range("A5").Value = "customer.Service@UOBgroup.com "
$ws.Range("B5").Value = "whitelisted"
$ws.Range("A5").Style = $ws.Range("A2").Style

# --- Misc sheet/page tweaks that came along with this resave -------------
$ws.Range("K10").Select()
$ws.PageSetup.Orientation = 1
